$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that originally followed the
#    H1 title ("Play Fire of Egypt Free - Merkur Slot Review").
$d.Paragraphs.Item(2).Range.Delete()

# 2. Insert a new bold-styled paragraph ("Play Fire of Egypt Free - Merkur
#    Slot Review") right before the closing "Prompt: ..." paragraph at the
#    end of the document.
$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fire of Egypt Free - Merkur Slot Review</w:t></w:r></w:p><w:p/><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
[void]$insertionPoint.InsertXML($newParaXml)

# InsertXML leaves a spare empty paragraph behind (the trailing paragraph
# mark of the injected fragment merges into the following paragraph, not
# the preceding one) - drop it so only the new bold paragraph remains.
$spareIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($spareIndex).Range.Delete()

# 3. Swap out the old image-generation "Prompt: ..." text for the new
#    meta-description sentence, keeping the existing italic run formatting.
[void]$d.Content.Find.Execute(
    "Prompt: Create a cartoon style image featuring a happy Maya warrior with glasses for the game " + [char]34 + "Fire of Egypt" + [char]34 + ". The image should be vibrant and eye-catching, with the Maya warrior as the focal point. The warrior should be shown holding some type of ancient Egyptian artifact or symbol, such as a pharaoh scepter or a golden necklace, to tie in with the game's theme. The background of the image should feature some of the game's symbols, such as the head of Queen Nefertiti or Tutankhamun's golden sarcophagus. The overall style of the image should be fun and playful to appeal to players who enjoy online slot games.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Fire of Egypt, a Merkur slot game set in ancient Egypt. Play free and discover the mix of classic and modern gameplay mechanics.",
    2)
